$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Address list re-sorted by distance from the warehouse; rows 5-13 get new
# values and a new row 14 is appended.
$ws.Range("A5").Value = "11716 Oban Ave, Tampa, FL 33617, USA"
$ws.Range("A6").Value = "4011 E Busch Blvd, Tampa, FL 33617, USA"
$ws.Range("A7").Value = "9410 N 12th St, Tampa, FL 33612, USA"
$ws.Range("A8").Value = "11710 Phoenix Cir, Tampa, FL 33618, USA"
$ws.Range("A9").Value = "10928 Lynn Lake Cir, Tampa, FL 33625, USA"
$ws.Range("A10").Value = "11004 Lynn Lake Cir, Tampa, FL 33625, USA"
$ws.Range("A11").Value = "10932 Lynn Lake Cir, Tampa, FL 33625, USA"
$ws.Range("A12").Value = "11002 Lynn Lake Cir, Tampa, FL 33625, USA"
$ws.Range("A13").Value = "7201 W Linebaugh Ave, Tampa, FL 33625, USA"
$ws.Range("A14").Value = "4505 N Armenia Ave, Tampa, FL 33603, USA"

# Move the selection cursor to A2, matching the post-edit saved view state.
$ws.Range("A2").Select() | Out-Null
